$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New id (column B) / speaker_variant (column C) values per row, now generated
# with no levenshtein-distance canonicalisation (each raw id kept distinct).
$ws.Cells.Item(2, 2).Value = "#octavio"
$ws.Cells.Item(2, 3).Value = "Octavio"

$ws.Cells.Item(3, 2).Value = "#guyd"
$ws.Cells.Item(3, 3).Value = "Guyd"

$ws.Cells.Item(4, 2).Value = "#brane"
$ws.Cells.Item(4, 3).Value = "Brane"

$ws.Cells.Item(5, 2).Value = "#gu.-gal"
$ws.Cells.Item(5, 3).Value = "Gu. Gal"

$ws.Cells.Item(6, 2).Value = "#guydeon"
$ws.Cells.Item(6, 3).Value = "Guydeon"

$ws.Cells.Item(7, 2).Value = "#billinc"
$ws.Cells.Item(7, 3).Value = "Billinc"

$ws.Cells.Item(8, 2).Value = "#keyser"
$ws.Cells.Item(8, 3).Value = "Keyser"

$ws.Cells.Item(9, 2).Value = "#keyzer"
$ws.Cells.Item(9, 3).Value = "Keyzer"

$ws.Cells.Item(10, 2).Value = "#galdra"
$ws.Cells.Item(10, 3).Value = "Galdra"

$ws.Cells.Item(11, 2).Value = "#gald"
$ws.Cells.Item(11, 3).Value = "Gald"

$ws.Cells.Item(12, 2).Value = "#galda"
$ws.Cells.Item(12, 3).Value = "Galda"

$ws.Cells.Item(13, 2).Value = "#billinci"
$ws.Cells.Item(13, 3).Value = "Billinci"

$ws.Cells.Item(14, 2).Value = "#galdrad"
$ws.Cells.Item(14, 3).Value = "Galdrad"

$ws.Cells.Item(15, 2).Value = "#keyzer"
$ws.Cells.Item(15, 3).Value = "keyzer"

$ws.Cells.Item(16, 2).Value = "#keyser"
$ws.Cells.Item(16, 3).Value = "keyser"

$ws.Cells.Item(17, 2).Value = "#galra"
$ws.Cells.Item(17, 3).Value = "Galra"

$ws.Cells.Item(18, 2).Value = "#kluyzen"
$ws.Cells.Item(18, 3).Value = "Kluyzen"

$ws.Cells.Item(19, 2).Value = "#galdr"
$ws.Cells.Item(19, 3).Value = "Galdr"

$ws.Cells.Item(20, 2).Value = "#flavio"
$ws.Cells.Item(20, 3).Value = "Flavio"

$ws.Cells.Item(21, 2).Value = "#kluys"
$ws.Cells.Item(21, 3).Value = "Kluys"

$ws.Cells.Item(22, 2).Value = "#galdrade"
$ws.Cells.Item(22, 3).Value = "Galdrade"

$ws.Cells.Item(23, 2).Value = "#laura"
$ws.Cells.Item(23, 3).Value = "Laura"

$ws.Cells.Item(24, 2).Value = "#guyde"
$ws.Cells.Item(24, 3).Value = "Guyde"

# The export no longer marks a preferred variant per id (is_prefered column),
# so clear out every "x" mark that used to flag the preferred row (rows 2-12).
for ($r = 2; $r -le 12; $r++) {
  $ws.Cells.Item($r, 4).ClearContents()
}
